# Applies the "algorithm.docx" revision described by the commit:
#   - appends a clause to each of the five "Ask user for ..." bullets
#   - adds a new sub-bullet (list level 2 / ilvl=1) under each of those five
#     bullets describing the float conversion
#   - removes the four old high-level bullets about plugging variables into
#     equations / printing the result
#   - replaces them with the full worked-out algorithm (create variables,
#     branch on population change, output results)

$d = $word.ActiveDocument

function Find-Range([string]$text) {
    $rng = $d.Content
    [void]$rng.Find.Execute($text)
    return $rng
}

# Helper: given the text of an existing top-level (ilvl=0) bullet, append
# $extra text to that same bullet, then insert a brand-new sub-bullet
# (ilvl = 1) right after it containing $subText. Always anchors off a fresh
# Document.Content range (Paragraph.Range based ranges misbehave after
# InsertParagraphAfter in this host).
function Add-Clause-And-Subbullet([string]$anchorText, [string]$extra, [string]$subText) {
    $r = Find-Range($anchorText)
    $r.Collapse(0)
    $r.InsertAfter($extra)
    $r.Collapse(0)
    [void]$r.InsertParagraphAfter()
    $r.Collapse(0)
    [void]$r.MoveStart(1, 1)
    $r.InsertAfter($subText)
    $r.ListFormat.ListIndent()
}

Add-Clause-And-Subbullet "Ask user for time between births" " in seconds and set to birth rate" " convert input to float"
Add-Clause-And-Subbullet "Ask user for time between deaths" " in second and set to death rate" "Convert input to float"
Add-Clause-And-Subbullet "Ask user for time between migrations" " and set to immigrate rate" "Convert input to float"
Add-Clause-And-Subbullet "Ask user for current population" " and set to current population" "Convert input to float"
Add-Clause-And-Subbullet "Ask user for number of years in future projection" " and set to years into future" "Convert input to a float"

# Remove the four old bullets:
#   "plug variables into population change equations for total population change"
#   "plug in variables for future population equation"
#   "Plug in variables for future increase/decrease "
#   "Print if there was an increase, decrease, or no change in population"
# These are the last four paragraphs in the document body (just before sectPr).
$delStart = Find-Range("plug variables into population change equations for total population change")
$start = $delStart.Start
$bodyEnd = $d.Content.End
$delRange = $d.Range($start, $bodyEnd)
$delRange.Delete()

# Append the replacement bullets after the new "Convert input to a float" sub-bullet
# (now the last paragraph in the document). $script:curLevel tracks the ilvl of the
# most-recently-inserted bullet so each call only needs to Indent/Outdent by the
# (small, here always <= 1) delta to reach the requested level -- ListIndent /
# ListOutdent are relative, one-level-at-a-time operations.
$script:curLevel = 1

function Add-Bullet([int]$ilvl, [string]$text) {
    $r = $d.Content
    $r.Collapse(0)
    [void]$r.InsertParagraphAfter()
    $r.Collapse(0)
    $r.InsertAfter($text)
    while ($script:curLevel -lt $ilvl) {
        $r.ListFormat.ListIndent()
        $script:curLevel = $script:curLevel + 1
    }
    while ($script:curLevel -gt $ilvl) {
        $r.ListFormat.ListOutdent()
        $script:curLevel = $script:curLevel - 1
    }
}

Add-Bullet 0 "Create variable seconds per year"
Add-Bullet 1 "Set variable equal to 365 * 24 * 60 * 60"
Add-Bullet 0 "Create variable births per year"
Add-Bullet 1 "Set variable to seconds per year divided by birth rate"
Add-Bullet 0 "Create variable deaths per year"
Add-Bullet 1 "Set variable to seconds per year divided by death rate"
Add-Bullet 0 "Create variable immigrants per year"
Add-Bullet 1 "Set variable to seconds per year divided by immigrant rate"
Add-Bullet 0 "Create variable population change per year"
Add-Bullet 1 "Set variable to births per year plus immigrants per year minus deaths per year"
Add-Bullet 0 "Create variable future population "
Add-Bullet 1 "Set equal to current population plus population change per year times years into future"
Add-Bullet 0 "If future population is greater than current population"
Add-Bullet 1 "Set change to increased"
Add-Bullet 0 "If future population is less than current population"
Add-Bullet 1 "Set change to decreased"
Add-Bullet 0 "Otherwise set change to remained the same "

$lq = [char]0x2018
$rq = [char]0x2019
$output1 = "Output " + $lq + "{years into the future} years, the population will be {future population}" + $rq
$output2 = "Output " + $lq + "the population has {change}" + $rq
Add-Bullet 0 $output1
Add-Bullet 0 $output2
